$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.723.61'
$ws.Range("E2").Value = '  +0.40%  '

$ws.Range("D3").Value = '1.638.83'
$ws.Range("E3").Value = '  -0.08%  '

$ws.Range("E4").Value = '  +0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.57'
$ws.Range("E5").Value = '  +1.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.503'
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0623'
$ws.Range("E9").Value = '  -0.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.07'
$ws.Range("E10").Value = '  -0.01%  '

$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").Value = '1.869.32'
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").Value = '1.650.96'
$ws.Range("E13").Value = '  +0.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.15'
$ws.Range("E14").Value = '  -0.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.526'
$ws.Range("E15").Value = '  -0.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.61'
$ws.Range("E16").Value = '  -0.28%  '

$ws.Range("D17").Value = '26.716.14'
$ws.Range("E17").Value = '  +0.35%  '

$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  -1.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '212.96'
$ws.Range("E19").Value = '  -1.12%  '

$ws.Range("E20").Value = '  +0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.36'
$ws.Range("E21").Value = '  +0.42%  '

$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.21'
$ws.Range("E22").Value = '  -0.77%  '

$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.34'
$ws.Range("E23").Value = '  +5.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.24'
$ws.Range("E24").Value = '  -2.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.40'
$ws.Range("E25").Value = '  +0.26%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("E27").Value = '  -1.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.13'
$ws.Range("E28").Value = '  +0.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.65'
$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0506'
$ws.Range("E30").Value = '  -1.22%  '

$ws.Range("E31").Value = '  +1.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.39'
$ws.Range("E32").Value = '  +1.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.99'
$ws.Range("E33").Value = '  -0.30%  '

$ws.Range("D34").Value = '1.283.26'
$ws.Range("E34").Value = '  +0.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  -0.52%  '

$ws.Range("E36").Value = '  +0.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0176'
$ws.Range("E37").Value = '  -1.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.533'
$ws.Range("E38").Value = '  +0.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.812'
$ws.Range("E39").Value = '  -1.01%  '

$ws.Range("E40").Value = '  +0.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.803'
$ws.Range("E41").Value = '  -0.97%  '

$ws.Range("E42").Value = '  -1.23%  '

$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D43").Value = '1.778.97'
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.26'
$ws.Range("E44").Value = '  -2.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.93'
$ws.Range("E45").Value = '  +3.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.47'
$ws.Range("E46").Value = '  -0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.59'
$ws.Range("E47").Value = '  -0.82%  '

$ws.Range("E48").Value = '  +0.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.61'
$ws.Range("E49").Value = '  -1.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0961'
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("E51").Value = '  +0.18%  '
